# Scheduled-runner style refresh of the Excalibur Leve-flipping profit
# sheets: market-price driven columns (H..N) are rewritten per leve row
# with freshly pulled averages; only numeric value cells change, no
# formulas/formatting are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 219.46153
$ws.Range("J19").Value = 217.16667
$ws.Range("L19").Value = 217.16667
$ws.Range("N19").Value = -567.1666700000001
$ws.Range("H33").Value = 262.22726
$ws.Range("I33").Value = 235.26315
$ws.Range("J33").Value = 433
$ws.Range("K33").Value = 235.26315
$ws.Range("L33").Value = 433
$ws.Range("M33").Value = -6.263149999999996
$ws.Range("N33").Value = -891
$ws.Range("H112").Value = 1629.1224
$ws.Range("J112").Value = 1740.093
$ws.Range("L112").Value = 5220.279
$ws.Range("N112").Value = -7436.279
$ws.Range("H131").Value = 5984.85
$ws.Range("I131").Value = 5284.385
$ws.Range("J131").Value = 7285.7144
$ws.Range("K131").Value = 15853.155
$ws.Range("L131").Value = 21857.1432
$ws.Range("M131").Value = -10813.155
$ws.Range("N131").Value = -31937.1432
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 286415
$ws.Range("I11").Value = 667335.3
$ws.Range("J11").Value = 724.75
$ws.Range("K11").Value = 667335.3
$ws.Range("L11").Value = 724.75
$ws.Range("M11").Value = -667191.3
$ws.Range("N11").Value = -1012.75
$ws.Range("H25").Value = 1614.5714
$ws.Range("I25").Value = 625.75
$ws.Range("J25").Value = 2933
$ws.Range("K25").Value = 625.75
$ws.Range("L25").Value = 2933
$ws.Range("M25").Value = -223.75
$ws.Range("N25").Value = -3737
$ws.Range("H26").Value = 2625
$ws.Range("I26").Value = 1833.3334
$ws.Range("K26").Value = 1833.3334
$ws.Range("M26").Value = -1503.3334
$ws.Range("H32").Value = 5378764.5
$ws.Range("I32").Value = 5884718.5
$ws.Range("K32").Value = 5884718.5
$ws.Range("M32").Value = -5884431.5
$ws.Range("H61").Value = 2567611.2
$ws.Range("I61").Value = 2781395.5
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 2781395.5
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -2781183.5
$ws.Range("N61").Value = -2624
$ws.Range("H74").Value = 3127887.5
$ws.Range("I74").Value = 3789913.8
$ws.Range("J74").Value = 6905.5713
$ws.Range("K74").Value = 3789913.8
$ws.Range("L74").Value = 6905.5713
$ws.Range("M74").Value = -3789039.8
$ws.Range("N74").Value = -8653.5713
$ws.Range("H77").Value = 3127887.5
$ws.Range("I77").Value = 3789913.8
$ws.Range("J77").Value = 6905.5713
$ws.Range("K77").Value = 18949569
$ws.Range("L77").Value = 34527.85649999999
$ws.Range("M77").Value = -18945201
$ws.Range("N77").Value = -43263.85649999999
$ws.Range("H102").Value = 3578.6667
$ws.Range("I102").Value = 2773.5881
$ws.Range("J102").Value = 7000.25
$ws.Range("K102").Value = 2773.5881
$ws.Range("L102").Value = 7000.25
$ws.Range("M102").Value = -1151.5881
$ws.Range("N102").Value = -10244.25
$ws.Range("H132").Value = 1667164.8
$ws.Range("I132").Value = 2910935.2
$ws.Range("K132").Value = 8732805.600000001
$ws.Range("M132").Value = -8730275.600000001
$ws.Range("H136").Value = 2567611.2
$ws.Range("I136").Value = 2781395.5
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 8344186.5
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -8341636.5
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2138.5557
$ws.Range("I20").Value = 2257.7778
$ws.Range("K20").Value = 2257.7778
$ws.Range("M20").Value = -2010.7778
$ws.Range("H105").Value = 1494.8636
$ws.Range("I105").Value = 1544.45
$ws.Range("K105").Value = 1544.45
$ws.Range("M105").Value = 202.55
$ws.Range("H134").Value = 430653.1
$ws.Range("I134").Value = 518270.22
$ws.Range("J134").Value = 10090.8
$ws.Range("K134").Value = 1554810.66
$ws.Range("L134").Value = 30272.4
$ws.Range("M134").Value = -1552275.66
$ws.Range("N134").Value = -35342.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 826153.1
$ws.Range("I58").Value = 1123888.1
$ws.Range("K58").Value = 1123888.1
$ws.Range("M58").Value = -1123685.1
$ws.Range("H122").Value = 4632.154
$ws.Range("I122").Value = 2024.8
$ws.Range("J122").Value = 6261.75
$ws.Range("K122").Value = 6074.4
$ws.Range("L122").Value = 18785.25
$ws.Range("M122").Value = -3624.4
$ws.Range("N122").Value = -23685.25
$ws.Range("H132").Value = 5008326
$ws.Range("I132").Value = 9332.884
$ws.Range("J132").Value = 35716428
$ws.Range("K132").Value = 27998.652
$ws.Range("L132").Value = 107149284
$ws.Range("M132").Value = -25468.652
$ws.Range("N132").Value = -107154344
$ws.Range("H134").Value = 2853.3333
$ws.Range("I134").Value = 2897.5
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 8692.5
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -6157.5
$ws.Range("N134").Value = -12570
$ws.Range("H136").Value = 826153.1
$ws.Range("I136").Value = 1123888.1
$ws.Range("K136").Value = 3371664.3
$ws.Range("M136").Value = -3369114.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 534985.0600000001
$ws.Range("I122").Value = 498.14285
$ws.Range("J122").Value = 731901.3
$ws.Range("K122").Value = 4483.28565
$ws.Range("L122").Value = 6587111.7
$ws.Range("M122").Value = -2033.28565
$ws.Range("N122").Value = -6592011.7
$ws.Range("H132").Value = 3305.6
$ws.Range("I132").Value = 2198.5
$ws.Range("J132").Value = 4043.6667
$ws.Range("K132").Value = 19786.5
$ws.Range("L132").Value = 36393.0003
$ws.Range("M132").Value = -17256.5
$ws.Range("N132").Value = -41453.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1097316.9
$ws.Range("I132").Value = 1508074.6
$ws.Range("K132").Value = 4524223.800000001
$ws.Range("M132").Value = -4521693.800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2672.8
$ws.Range("I7").Value = 2361.1428
$ws.Range("J7").Value = 3400
$ws.Range("K7").Value = 2361.1428
$ws.Range("L7").Value = 3400
$ws.Range("M7").Value = -2249.1428
$ws.Range("N7").Value = -3624
$ws.Range("H122").Value = 3756.3057
$ws.Range("I122").Value = 3563
$ws.Range("J122").Value = 4098.3076
$ws.Range("K122").Value = 10689
$ws.Range("L122").Value = 12294.9228
$ws.Range("M122").Value = -8239
$ws.Range("N122").Value = -17194.9228
$ws.Range("H126").Value = 2672.8
$ws.Range("I126").Value = 2361.1428
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 7083.428400000001
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -4613.428400000001
$ws.Range("N126").Value = -15140
$ws.Range("H132").Value = 2163266
$ws.Range("I132").Value = 2471832
$ws.Range("J132").Value = 3304.5
$ws.Range("K132").Value = 7415496
$ws.Range("L132").Value = 9913.5
$ws.Range("M132").Value = -7412966
$ws.Range("N132").Value = -14973.5
$ws.Range("H136").Value = 6359.8667
$ws.Range("I136").Value = 4439.9
$ws.Range("K136").Value = 13319.7
$ws.Range("M136").Value = -10769.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6291867
$ws.Range("I132").Value = 6942173.5
$ws.Range("K132").Value = 20826520.5
$ws.Range("M132").Value = -20823990.5
$ws.Range("H136").Value = 7312268
$ws.Range("I136").Value = 8448888
$ws.Range("J136").Value = 5428.5713
$ws.Range("K136").Value = 25346664
$ws.Range("L136").Value = 16285.7139
$ws.Range("M136").Value = -25344114
$ws.Range("N136").Value = -21385.7139
